$wb = $excel.ActiveWorkbook

$uk = $wb.Worksheets.Item("UK")

# --- Romania sheet: clone UK (keeps exact column widths/styles/merges),
#     strip the FC32AR/FC32DR repeater rows, relabel header row ---
[void]$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"

[void]$romania.Rows.Item(8).Delete()
[void]$romania.Rows.Item(8).Delete()

$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3537/T3551"
[void]$romania.Range("B4").Select()

# --- Slovakia sheet: same clone-and-trim recipe ---
[void]$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

[void]$slovakia.Rows.Item(8).Delete()
[void]$slovakia.Rows.Item(8).Delete()

$slovakia.Range("B2").Value = "Slovakia market"
$slovakia.Range("B4").Value = "NGC-4306/T3564/T3576"
$slovakia.Range("B4").Style = "Normal"
[void]$slovakia.Range("B4").Select()

# UK itself ends up with a "select-all" selection and loses the active tab
[void]$uk.Cells.Select()

# Slovakia (last sheet) is the active tab on save
[void]$slovakia.Activate()
